$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 136 (Indice 135) ----
$ws.Range("A136").Value = 135
$ws.Range("B136").Value = "romania"
$ws.Range("C136").Value = "liga-1"
$ws.Range("D136").Value = "2023-2024"
$ws.Range("E136").Value = 45257.6875
$ws.Range("F136").Value = "U Craiova 1948"
$ws.Range("G136").Value = 2
$ws.Range("H136").Value = "Sepsi Sf. Gheorghe"
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = 2.7
$ws.Range("K136").Value = "20/11/2023 16:42"
$ws.Range("L136").Value = 2.88
$ws.Range("M136").Value = "27/11/2023 16:25"
$ws.Range("N136").Value = 3.14
$ws.Range("O136").Value = "20/11/2023 16:42"
$ws.Range("P136").Value = 3.09
$ws.Range("Q136").Value = "27/11/2023 16:25"
$ws.Range("R136").Value = 2.74
$ws.Range("S136").Value = "20/11/2023 16:42"
$ws.Range("T136").Value = 2.68
$ws.Range("U136").Value = "27/11/2023 16:25"
$ws.Range("V136").Value = "https://www.betexplorer.com/football/romania/liga-1/fc-u-craiova-sepsi/hE5UCxHc/"

# ---- Row 137 (Indice 136) ----
$ws.Range("A137").Value = 136
$ws.Range("B137").Value = "romania"
$ws.Range("C137").Value = "liga-1"
$ws.Range("D137").Value = "2023-2024"
$ws.Range("E137").Value = 45257.8125
$ws.Range("F137").Value = "FC Rapid Bucuresti"
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = "U. Cluj"
$ws.Range("I137").Value = 3
$ws.Range("J137").Value = 1.72
$ws.Range("K137").Value = "23/11/2023 11:42"
$ws.Range("L137").Value = 1.9
$ws.Range("M137").Value = "27/11/2023 19:27"
$ws.Range("N137").Value = 3.7
$ws.Range("O137").Value = "23/11/2023 11:42"
$ws.Range("P137").Value = 3.49
$ws.Range("Q137").Value = "27/11/2023 19:26"
$ws.Range("R137").Value = 4.86
$ws.Range("S137").Value = "23/11/2023 11:42"
$ws.Range("T137").Value = 4.26
$ws.Range("U137").Value = "27/11/2023 19:25"
$ws.Range("V137").Value = "https://www.betexplorer.com/football/romania/liga-1/rapid-bucuresti-universitatea-cluj/Ak0ZBdW3/"

# ---- Match formatting of the preceding data row (row 135) ----
# Column A uses the bold/bordered/centered header-like style; column E uses
# the custom date-time number format. Copy formats only (xlPasteFormats)
# so no duplicate style entries get created.
$ws.Range("A135").Copy() | Out-Null
$ws.Range("A136:A137").PasteSpecial(-4122) | Out-Null

$ws.Range("E135").Copy() | Out-Null
$ws.Range("E136:E137").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
